$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# StedenInformatie (sheet2): restructure the table
# ---------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("StedenInformatie")

# Insert a new column A (shifts: old A->B city, B->C inwoners,
# C->D testees-formula, D->E treshold-formula). Column A's old
# width (16) rides along to the new column B, matching the diff.
$ws2.Columns.Item(1).Insert() | Out-Null

# New column A: running index 1..12
for ($i = 1; $i -le 12; $i++) {
    $ws2.Cells.Item($i + 1, 1).Value = $i
}

# The old "treshold" formula column (now E) becomes a flat constant 10
for ($r = 2; $r -le 13; $r++) {
    $ws2.Cells.Item($r, 5).Value = 10
}

# New totals row
$ws2.Range("D14").Formula = "=SUM(D2:D13)"

# Headers (string-table order matters: Stad_naam before LL)
$ws2.Range("B1").Value = "Stad_naam"
$ws2.Range("A1").Value = "LL"

# A1 header style: bold + right aligned
$ws2.Range("A1").Font.Bold = $true
$ws2.Range("A1").HorizontalAlignment = -4152   # xlRight

# Column widths: col B (was col A) 16 chars, col E ~23.875 chars
$ws2.Columns.Item(5).ColumnWidth = 23

# ---------------------------------------------------------------
# RandomPopulation (sheet1): selection moves, no longer tabSelected
# ---------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("RandomPopulation")
$ws1.Range("A9:A25").Select() | Out-Null

# Sheet2 becomes the active/selected sheet (do this LAST so it sticks)
$ws2.Activate() | Out-Null
$ws2.Range("G4").Select() | Out-Null

Write-Host "done"
